# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.599.92"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.800.10"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'227.39"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'32.84"
$ws.Range("E8").Value = "  +3.83%  "
$ws.Range("D9").Value = "'0.296"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.0950"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "2.061.87"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.15"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.798.13"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "34.603.78"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "'4.33"
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").Value = "'68.95"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").Value = "0.0₃0803"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'247.29"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").Value = "'172.17"
$ws.Range("E24").Value = "  +5.99%  "
$ws.Range("D25").Value = "'2.06"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'7.33"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").Value = "'16.62"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'4.05"
$ws.Range("E30").Value = "  +9.23%  "
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "'3.81"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("D35").Value = "1.433.42"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").Value = "  +6.79%  "
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "'84.80"
$ws.Range("E40").Value = "  +5.91%  "
$ws.Range("E41").Value = "  +3.03%  "
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").Value = "'2.76"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").Value = "'13.82"
$ws.Range("E44").Value = "  +2.82%  "
$ws.Range("D45").Value = "'0.0527"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").Value = "'6.11"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "1.961.78"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").Value = "'105.35"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("E51").Value = "  -4.99%  "
